$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion message text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.55 = 18090.95 pesos`n✅ 18090.95 pesos = 4.53 = 967.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 220
$ws2.Range("O10").Value = 3980.01
$ws2.Range("N12").Value = 3992
$ws2.Range("O12").Value = 213.5
